try {
  $word.OrganizerCopy("a","b","c","d")
  Write-Output "called, no exception"
} catch {
  Write-Output "EXC: $_"
}
